# Append the latest Nalco ingot PDF run-log entry (row 36) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count + 1, 1).Row
if ($newRow -lt 36) { $newRow = 36 }

$ws.Cells.Item($newRow, 1).Value = "2025-08-20 09:40:14 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-20 15:10:14 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

# Match the formatting used by the existing rows (style index 3: centered
# text, default font/border), copied from the row directly above.
$srcRange = $ws.Range($ws.Cells.Item($newRow - 1, 1), $ws.Cells.Item($newRow - 1, 8))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 8))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)
